# Update "想去人数" (want-to-go count) figures to reflect a fresh scrape
# as described by commit "Update gh-pages to output generated at 456a3b4".
#
# Sheet "展览" (exhibitions):
#   row 2  (ACGN夏日游园会第七回-泳池派对)        F: 618  -> 619
#   row 10 (第十五届次元之门动漫游戏博览会)        F: 4985 -> 4991
#   row 11 (首届AT次元时代动漫游戏嘉年华)          F: 4686 -> 4692
#
# Sheet "演出" (performances):
#   row 2  (四月是你的谎言 钢琴小提琴音乐集)       F: 68   -> 69
#
# Sheet "全部类型" (all types - aggregate of every other sheet):
#   row 2  (ACGN夏日游园会第七回-泳池派对)          F: 618  -> 619
#   row 10 (第十五届次元之门动漫游戏博览会)         F: 4985 -> 4991
#   row 11 (首届AT次元时代动漫游戏嘉年华)           F: 4686 -> 4692
#   row 17 (四月是你的谎言 钢琴小提琴音乐集)        F: 68   -> 69

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 619
$wsExhibit.Range("F10").Value = 4991
$wsExhibit.Range("F11").Value = 4692

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 69

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 619
$wsAll.Range("F10").Value = 4991
$wsAll.Range("F11").Value = 4692
$wsAll.Range("F17").Value = 69
